# Weekly fruit/vegetable price update: insert this week's new price
# observation as a new row 278 ("Feria Lagunitas de Puerto Montt" /
# Ciboulette), pushing all subsequent rows (old 278..388) down by one
# (new 279..389).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 278; Excel shifts rows 278-388 down
# to 279-389 and extends the used range (dimension) automatically.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row 278 with the new weekly record.
$ws.Cells.Item(278, 1).Value  = 4
$ws.Cells.Item(278, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(278, 3).Value  = "Los Lagos"
$ws.Cells.Item(278, 4).Value  = 45141
$ws.Cells.Item(278, 5).Value  = 10
$ws.Cells.Item(278, 6).Value  = 100112039
$ws.Cells.Item(278, 7).Value  = "Ciboulette"
$ws.Cells.Item(278, 8).Value  = "Sin especificar"
$ws.Cells.Item(278, 9).Value  = "Primera"
$ws.Cells.Item(278, 10).Value = 80
$ws.Cells.Item(278, 11).Value = 2500
$ws.Cells.Item(278, 12).Value = 3500
$ws.Cells.Item(278, 13).Value = 3000
$ws.Cells.Item(278, 14).Value = "`$/docena de atados"
$ws.Cells.Item(278, 15).Value = "Región Metropolitana"
$ws.Cells.Item(278, 16).Value = 1000
$ws.Cells.Item(278, 17).Value = 3
$ws.Cells.Item(278, 18).Value = "Hortaliza"
